$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Bill"
$ws.Range("B1").Value = "Admin"
$ws.Range("C1").Value = "Stockholm"
$ws.Range("D1").Value = "Audi"

$ws.Range("A2").Value = "George"
$ws.Range("B2").Value = "User"
$ws.Range("C2").Value = "Gothenburg"
$ws.Range("D2").Value = "Volvo"
